# Insert a new data row at row 205 (pushing former rows 205-265 down to
# 206-266) and populate it with a new "Plátano" price observation for
# Agrícola del Norte S.A. de Arica, dated 2022-07-13 (Excel serial 44755).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 205..265 down one position, creating a blank row 205
# (the new row inherits formatting, incl. the date number format on col D,
# from the row above it, same as native Excel "Insert").
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record.
$ws.Cells.Item(205,1).Value  = 1
$ws.Cells.Item(205,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(205,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(205,4).Value  = 44755
$ws.Cells.Item(205,5).Value  = 15
$ws.Cells.Item(205,6).Value  = "Fruta"
$ws.Cells.Item(205,7).Value  = 100108
$ws.Cells.Item(205,8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(205,9).Value  = 100108006
$ws.Cells.Item(205,10).Value = "Plátano"
$ws.Cells.Item(205,11).Value = "Sin especificar"
$ws.Cells.Item(205,12).Value = "Pintón"
$ws.Cells.Item(205,13).Value = 120
$ws.Cells.Item(205,14).Value = 25000
$ws.Cells.Item(205,15).Value = 26000
$ws.Cells.Item(205,16).Value = 25500
$ws.Cells.Item(205,17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(205,18).Value = "Ecuador"
$ws.Cells.Item(205,19).Value = 1275
$ws.Cells.Item(205,20).Value = 20
